$d = $word.ActiveDocument

function Get-ParaEndRange($para) {
    # Returns a collapsed range positioned right before the paragraph mark,
    # i.e. the correct insertion point for appending text at the end of a paragraph.
    $r = $para.Range
    return $d.Range($r.End - 1, $r.End - 1)
}

# ---------------------------------------------------------------------------
# 1) "Mozliwosc obracania pliku PDF" + "." -> merge into a single run ending
#    with a period: "Mozliwosc obracania pliku PDF."
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("Możliwość obracania pliku PDF.", $true, $false, $false, $false, $false, `
    $true, 1, $false, "Możliwość obracania pliku PDF.", 2) | Out-Null

# ---------------------------------------------------------------------------
# Locate the paragraphs we need to edit by their (stable) leading text.
# ---------------------------------------------------------------------------
$paraDzial = $null
$paraNrDok = $null
$paraSortowanie = $null
$paraRodzaj = $null
$paraZmeczenie = $null

foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t.StartsWith("dodać pole „Dział teczki”")) {
        $paraDzial = $p
    } elseif ($t.StartsWith("dodać pole ")) {
        $paraNrDok = $p
    } elseif ($t.StartsWith("sortowanie tabeli po numerze dokumentu")) {
        $paraSortowanie = $p
    } elseif ($t.StartsWith("w kolumnie ") -and $t.Contains("mnemonik")) {
        $paraRodzaj = $p
    } elseif ($t.Contains("zmęczenie pracownika opisującego")) {
        $paraZmeczenie = $p
    }
}

# ---------------------------------------------------------------------------
# 2) "dodać pole „Dział teczki” " -> append " " + red "//Michał/W BE zrobione, pozostał FE"
#    (written as four separate red runs, matching the original authoring)
# ---------------------------------------------------------------------------
$ip = Get-ParaEndRange $paraDzial
$ip.InsertAfter(" ")
$ip.Collapse(0)

$ip.InsertAfter("//")
$ip.Font.Color = 255
$ip.Collapse(0)

$ip.InsertAfter("Michał/")
$ip.Font.Color = 255
$ip.Collapse(0)

$ip.InsertAfter("W BE ")
$ip.Font.Color = 255
$ip.Collapse(0)

$ip.InsertAfter("zrobione, pozostał FE")
$ip.Font.Color = 255
$ip.Collapse(0)

# ---------------------------------------------------------------------------
# 3) "dodać pole ... (1A, 2A, 1B, 2B…) " -> append red "//Michał/W BE zrobione, pozostał FE"
# ---------------------------------------------------------------------------
$ip = Get-ParaEndRange $paraNrDok
$ip.InsertAfter("//Michał/W BE zrobione, pozostał FE")
$ip.Font.Color = 255

# ---------------------------------------------------------------------------
# 4) "sortowanie tabeli po numerze dokumentu." -> append " " + red annotation
# ---------------------------------------------------------------------------
$ip = Get-ParaEndRange $paraSortowanie
$ip.InsertAfter(" ")
$ip.Collapse(0)

$ip.InsertAfter("//Michał/W BE zrobione – sortowanie podwójne, w pierwszej kolejności po dziale teczki, potem po numerze dokumentu. Pozostał FE.")
$ip.Font.Color = 255

# ---------------------------------------------------------------------------
# 5) "w kolumnie Rodzaj dokumentu pełna nazwa, a nie mnemonik." -> add _GoBack
#    bookmark at the end of this paragraph (this also removes it from its old
#    location further below, exactly like real Word moves a single _GoBack
#    bookmark on every edit).
# ---------------------------------------------------------------------------
$ip = Get-ParaEndRange $paraRodzaj
$d.Bookmarks.Add("_GoBack", $ip) | Out-Null

# ---------------------------------------------------------------------------
# 6) Merge the two runs that used to sit around the old _GoBack bookmark into
#    a single run: " zakladac, biorac pod uwage zmeczenie pracownika
#    opisujacego i monotonnosc czynnosci."
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("zakładać, biorąc pod uwagę zmęczenie pracownika opisującego i monotonność czynności.", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "zakładać, biorąc pod uwagę zmęczenie pracownika opisującego i monotonność czynności.", 2) | Out-Null
